$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The logger dropped the oldest reading (row 6) and rows 2-5 now hold a newer
# batch of junction-flooding samples (custom-accuracy export, "1000 data points" run).
$ws.Range("A6:AH6").EntireRow.Delete()

# Row 2: Time + J1..J33 readings
$rowValues = @(45153.50694444445, 12.343, 8.237, 3.467, 26.724, 19.801, 9.406000000000001, 27.816, 15.125, 5.951, 8.669, 10.525, 11.375, 3.134, 9.775, 13.324, 8.832000000000001, 2.74, 1.542, 141.399, 27.009, 9.023, 17.239, 8.882, 2.65, 15.251, 7.97, 7.4, 8.521000000000001, 11.028, 2.858, 25.233, 4.735, 11.279)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
}

# Row 3: Time + J1..J33 readings
$rowValues = @(45153.51388888889, 19.096, 13.904, 1.849, 41.735, 33.303, 14.864, 55.68, 23.269, 10.095, 14.766, 16.689, 17.835, 4.83, 15.038, 21.202, 12.977, 1.393, 1.083, 221.58, 42.024, 13.881, 27.891, 14.577, 2.54, 27.904, 12.261, 11.022, 12.917, 17.495, 1.162, 50.999, 7.663, 17.355)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(3, $col).Value = $rowValues[$col - 1]
}

# Row 4: Time + J1..J33 readings
$rowValues = @(45153.52083333334, 0.355, 0.091, 0.765, 0.848, 0, 0, 9.651999999999999, 0.582, 0.143, 0.244, 0.193, 0.136, 0, 0.376, 0.639, 0.594, 0.82, 0.336, 0, 1.754, 0.347, 1.107, 0.386, 0.425, 4.308, 0.291, 0.468, 0.485, 0.363, 0.669, 9.766999999999999, 0.053, 0.448)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(4, $col).Value = $rowValues[$col - 1]
}

# Row 5: Time + J1..J33 readings
$rowValues = @(45153.52777777778, 12.87, 9.550000000000001, 0.98, 28.09, 22.66, 10.36, 35.42, 15.71, 6.75, 10.41, 11.17, 11.75, 3.16, 10.15, 14.31, 8.640000000000001, 0.73, 0.61, 147.59, 28.18, 9.369999999999999, 18.77, 9.890000000000001, 1.62, 17.56, 8.26, 7.4, 8.68, 11.86, 0.5, 31.9, 5.24, 11.7)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(5, $col).Value = $rowValues[$col - 1]
}

# Column widths widened slightly (character units -> stored width = chars + 5/6)
$ws.Range("B1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("C1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("F1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("G1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("I1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("J1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("K1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("L1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("M1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("O1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("P1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("T1").EntireColumn.ColumnWidth = 8.166666666666666
$ws.Range("V1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("W1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("X1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Z1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AA1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AB1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AC1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AD1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AF1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AH1").EntireColumn.ColumnWidth = 7.166666666666667
